$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain TEXT (inline strings
# in the source file), so force Text number format before assigning, exactly
# as a human would do in Excel to stop auto-conversion to a real number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.32"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.74"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.658"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05823"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.410"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.321"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07647"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03236"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02928"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09241"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001662"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.316"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04744"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005990"
$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006251"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005473"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001069"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.178"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3329"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1224"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0009998"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04293"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007113"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003602"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009527"

$ws.Range("E45").Value = "44ACDXExchangeACXT"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005447"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7853"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1025"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01011"
